# "Taking the latest code" - refresh the "Test Cases" sheet with the
# latest set of authored test cases: 3 new rows (TestCase_A21-A23), and
# flip TestCase_A16's (row 17) result from SKIP to PASS now that it has
# been re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# TestCase_A16 (row 17) now passes.
$ws.Cells.Item(17, 5).Value = "PASS"

# New row 22 - TestCase_A21
$ws.Cells.Item(22, 1).Value = "TestCase_A21"
$ws.Cells.Item(22, 2).Value = "OPQA-399"
$ws.Cells.Item(22, 3).Value = "Verify View additional email preferences link is working"
$ws.Cells.Item(22, 4).Value = "Y"
$ws.Cells.Item(22, 5).Value = "SKIP"

# New row 23 - TestCase_A22
$ws.Cells.Item(23, 1).Value = "TestCase_A22"
$ws.Cells.Item(23, 2).Value = "OPQA-854,OPQA-853"
$ws.Cells.Item(23, 3).Value = 'Verify that the  checkbox  is present and can be modified for "Receive email notifications for likes,comments and other activity" is working correctly.'
$ws.Cells.Item(23, 4).Value = "Y"
$ws.Cells.Item(23, 5).Value = "SKIP"
$ws.Rows.Item(23).RowHeight = 30

# New row 24 - TestCase_A23
$ws.Cells.Item(24, 1).Value = "TestCase_A23"
$ws.Cells.Item(24, 2).Value = "OPQA-527"
$ws.Cells.Item(24, 3).Value = "Verify change password link in the account page is working correctly."
$ws.Cells.Item(24, 4).Value = "Y"
$ws.Cells.Item(24, 5).Value = "SKIP"

# Match the borders used by the rest of the table for the 3 new rows.
$ws.Range("A22:E24").Borders.LineStyle = 1
$ws.Range("C22:C23").WrapText = $true

# Leave the cursor where the author left it while reviewing the new rows.
$ws.Activate()
$ws.Range("D17").Select()

$ws.PageSetup.Orientation = 1
